# Apply the "add on feature test plan" edits to AddOn_TestData.xlsx
#
# Summary of intended changes (see commit diff):
#   - shared text "AddOn Test Suite" (cell C1 on every sheet) -> "Web Data 69"
#   - AddProductCategory1!G1/G2 numeric values swapped (10 <-> 0.25)
#   - AddProductCategory1 rows 1-2 row height normalised to 15
#   - Active / selected worksheet moves from "AddCustomer" to "AddOrder"
#     (workbook tabRatio/activeTab + per-sheet tabSelected + top-left cell)

$wb = $excel.ActiveWorkbook

$wsProductCategory = $wb.Worksheets.Item("AddProductCategory1")
$wsCustomer        = $wb.Worksheets.Item("AddCustomer")
$wsOrder           = $wb.Worksheets.Item("AddOrder")

# --- Shared text update: "AddOn Test Suite" -> "Web Data 69" --------------
# The string lives in the shared strings table and is referenced from cell
# C1 on every sheet, so update it everywhere it appears.
$wsProductCategory.Range("C1").Value = "Web Data 69"
$wsCustomer.Range("C1").Value = "Web Data 69"
$wsOrder.Range("C1").Value = "Web Data 69"

# --- Swap the two numeric values on AddProductCategory1 -------------------
$wsProductCategory.Range("G1").Value = 0.25
$wsProductCategory.Range("G2").Value = 10

# --- Normalise row heights for the edited rows -----------------------------
$wsProductCategory.Rows.Item(1).RowHeight = 15
$wsProductCategory.Rows.Item(2).RowHeight = 15

# --- Move the active/selected sheet from AddCustomer to AddOrder ----------
# (this also moves workbook activeTab to index 2 and sets tabSelected on the
# AddOrder sheet view while clearing it elsewhere). The current selected
# cell (C1) is left untouched - only the view's top-left (scroll) cell moves
# back to A1.
$wsOrder.Activate()

$win = $excel.ActiveWindow
$win.ScrollColumn = 1
$win.ScrollRow = 1
$win.TabRatio = 993
